$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- R11 gets a "See below" note pointing at the new alternatives table ---
# (written first so "See below" becomes the first newly-introduced shared string)
$ws.Range("H42").Value2 = "See below"

# --- New section below the totals row: "R11 other options" ---
$ws.Range("A69").Value2 = "R11 other options"

# --- IC block shifts up (IC1 removed, IC2/IC3/IC4 shift, new IC5 added) ---
# Row 56 was IC1 -> becomes IC2 (reuses existing shared string "IC2")
$ws.Range("A56").Value2 = "IC2"

# Row 57 was IC2 -> becomes IC3, gets a new "No" flag in column G (both reused strings)
$ws.Range("A57").Value2 = "IC3"
$ws.Range("G57").Value2 = "No"

# Row 58 was IC3 -> becomes IC4 (unchanged otherwise, reused string)
$ws.Range("A58").Value2 = "IC4"

# Row 59 was IC4 -> becomes IC5 (new shared string), gets a new "No" flag in column G (reused)
$ws.Range("A59").Value2 = "IC5"
$ws.Range("G59").Value2 = "No"

# --- Finish populating the "R11 other options" alternatives table ---
$ws.Range("A73").Value2 = "0 Ohm resistor"
$ws.Range("B73").Value2 = 0
$ws.Range("C73").Value2 = "Thick Film"
$ws.Range("D73").Value2 = "Mouser"
$ws.Range("E73").Value2 = "652-CR1206-J/-000ELF"
$ws.Range("F73").Value2 = 0.1

$ws.Range("B70").Value2 = 200
$ws.Range("C70").Value2 = "Thick Film"
$ws.Range("D70").Value2 = "Mouser"
$ws.Range("E70").Value2 = "652-CR1206FX-2000ELF"
$ws.Range("F70").Value2 = 0.1

$ws.Range("B69").Value2 = 500
$ws.Range("C69").Value2 = "Thin Film"
$ws.Range("D69").Value2 = "Mouser"
$ws.Range("E69").Value2 = "71-PTN1206Y5000BST1"
$ws.Range("F69").Value2 = 1.04

$ws.Range("B71").Value2 = 100
$ws.Range("C71").Value2 = "Thick Film"
$ws.Range("D71").Value2 = "Mouser"
$ws.Range("E71").Value2 = "71-RCC1206100RJNEA"
$ws.Range("F71").Value2 = 0.27

# --- Update the saved cursor position / selection to match the end state ---
$ws.Range("G78").Select()
